## Add a new "title" column before the existing "Resultat" column (col F),
## shifting the old F/G... columns one place to the right, and select F8.
## This mirrors: "I save the value of '(.*)-(.*)' in '(.*)' column of data
## output provider." -> a new data-output-provider "title" column is added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at F; existing F (and everything to its right)
# shifts right by one column (F->G, G->H, ...). Excel copies the formatting
# of the column to the left (E) into the freshly inserted column.
$ws.Columns("F:F").Insert()

# The new column header (old G1, which now holds "Resultat") keeps its
# original header style; copy that style onto the new F1 header cell before
# giving it its own text so the header row formatting stays consistent.
$ws.Range("G1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F1").Value = "title"

# Match the new column's width to its neighbouring data column (E).
$ws.Columns("F:F").ColumnWidth = 14.584

# Restore the active selection to reflect where the user finished editing.
$ws.Range("F8").Select()
